{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the 25 two-digit division prompts in the practice-sheet table,\n// in document order, with their new values per the commit diff.\n\nconst replacements = [\n  \"23\u00f76=\", \"67\u00f79=\",\n  \"79\u00f78=\", \"62\u00f77=\",\n  \"58\u00f75=\", \"82\u00f74=\",\n  \"83\u00f78=\", \"82\u00f73=\",\n  \"96\u00f73=\", \"10\u00f76=\",\n  \"53\u00f79=\", \"12\u00f73=\",\n  \"11\u00f78=\", \"29\u00f78=\",\n  \"83\u00f72=\", \"26\u00f79=\",\n  \"50\u00f75=\", \"67\u00f78=\",\n  \"82\u00f76=\", \"69\u00f75=\",\n  \"61\u00f73=\", \"14\u00f78=\",\n  \"86\u00f76=\", \"58\u00f79=\",\n  \"89\u00f79=\", \"60\u00f73=\",\n  \"80\u00f73=\", \"85\u00f74=\",\n  \"25\u00f77=\", \"95\u00f77=\",\n  \"96\u00f73=\", \"90\u00f73=\",\n  \"88\u00f73=\", \"25\u00f76=\",\n  \"27\u00f74=\", \"31\u00f79=\",\n  \"80\u00f74=\", \"35\u00f79=\",\n  \"52\u00f73=\", \"46\u00f79=\",\n  \"58\u00f78=\", \"64\u00f79=\",\n  \"67\u00f77=\", \"24\u00f74=\",\n  \"16\u00f75=\", \"26\u00f72=\",\n  \"80\u00f74=\", \"10\u00f73=\",\n  \"52\u00f78=\", \"97\u00f77=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Walk the document's paragraphs in order, and for each one whose text\n// exactly matches the *next* expected \"old\" value in our ordered list,\n// replace it with the corresponding \"new\" value. This mirrors the\n// sequential, position-based nature of the diff (some old values like\n// \"96\u00f73=\" and \"80\u00f74=\" repeat, and each occurrence maps to a different\n// replacement based on its position in the document).\nlet pos = 0;\nfor (let i = 0; i < paragraphs.items.length && pos < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const oldText = replacements[pos];\n  const newText = replacements[pos + 1];\n  if (para.text === oldText) {\n    para.insertText(newText, \"Replace\");\n    pos += 2;\n  }\n}\n\nawait context.sync();\n\nif (pos !== replacements.length) {\n  throw new Error(\n    `Only matched ${pos / 2} of ${replacements.length / 2} expected cells.`\n  );\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces the 25 two-digit division prompts in the practice-sheet table,\n# in document order, with their new values per the commit diff.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"23\u00f76=\", \"67\u00f79=\",\n    \"79\u00f78=\", \"62\u00f77=\",\n    \"58\u00f75=\", \"82\u00f74=\",\n    \"83\u00f78=\", \"82\u00f73=\",\n    \"96\u00f73=\", \"10\u00f76=\",\n    \"53\u00f79=\", \"12\u00f73=\",\n    \"11\u00f78=\", \"29\u00f78=\",\n    \"83\u00f72=\", \"26\u00f79=\",\n    \"50\u00f75=\", \"67\u00f78=\",\n    \"82\u00f76=\", \"69\u00f75=\",\n    \"61\u00f73=\", \"14\u00f78=\",\n    \"86\u00f76=\", \"58\u00f79=\",\n    \"89\u00f79=\", \"60\u00f73=\",\n    \"80\u00f73=\", \"85\u00f74=\",\n    \"25\u00f77=\", \"95\u00f77=\",\n    \"96\u00f73=\", \"90\u00f73=\",\n    \"88\u00f73=\", \"25\u00f76=\",\n    \"27\u00f74=\", \"31\u00f79=\",\n    \"80\u00f74=\", \"35\u00f79=\",\n    \"52\u00f73=\", \"46\u00f79=\",\n    \"58\u00f78=\", \"64\u00f79=\",\n    \"67\u00f77=\", \"24\u00f74=\",\n    \"16\u00f75=\", \"26\u00f72=\",\n    \"80\u00f74=\", \"10\u00f73=\",\n    \"52\u00f78=\", \"97\u00f77=\"\n)\n\n# Walk the document's paragraphs in order, and for each one whose (trimmed)\n# text exactly matches the *next* expected \"old\" value in our ordered list,\n# replace it with the corresponding \"new\" value. This mirrors the\n# sequential, position-based nature of the diff (some old values like\n# \"96\u00f73=\" and \"80\u00f74=\" repeat, and each occurrence maps to a different\n# replacement based on its position in the document).\n$pos = 0\nforeach ($p in $d.Paragraphs) {\n    if ($pos -ge $replacements.Count) { break }\n    $r = $p.Range\n    $clean = $r.Text.TrimEnd([char]13, [char]7)\n    $oldText = $replacements[$pos]\n    $newText = $replacements[$pos + 1]\n    if ($clean -eq $oldText) {\n        $r.Text = $newText\n        $pos = $pos + 2\n    }\n}\n\nif ($pos -ne $replacements.Count) {\n    throw \"Only matched $($pos / 2) of $($replacements.Count / 2) expected cells.\"\n}\n"}
